# "Add more work log functions"
#
# Layout of the "Varpas 1" worksheet: row 1 is a header (day numbers 1..30,
# then "Kopā" (Total) in column AF and "Dienas" (Days) in column AG). Rows
# 2-4 hold one worked-hours log per person, column A being their name and
# columns B..AE the hours logged for each day of the month.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 used to belong to "Signe Zalužinska"; it is now a new person's log.
$ws.Range("A3").Value = "Jauns efwefew"

# Row 4 used to belong to "fwefwf wfwefwe"; it is now Signe Zalužinska's log.
$ws.Range("A4").Value = "Signe Zalužinska"

# Row 3 (now "Jauns efwefew"): day 10 (column K) hours reset to 0, which
# zeroes out the row's Total (AF) and Days (AG = Total / 24) as well.
$ws.Range("K3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0

# Row 4 (now "Signe Zalužinska"): day 9 (column J) becomes an empty/"0"
# text entry instead of a numeric 0, and day 10 (column K) changes from
# 7 to 9 hours, bumping the Total (AF) and Days (AG = Total / 24).
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "0"
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Value = 9
$ws.Range("AF4").Value = 9
$ws.Range("AG4").Value = 0.375
